# Generate Report for Handback
#
# The handback for ebb52c47-83da-4006-ae86-ec3795f654d8 (the second file in
# the report) has progressed: it was handed off again and a new handback
# was received, so refresh the tracked timestamps across all three sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# "Latest HO Xliff Generate Date" for the ebb52c47 row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-22 20:47:46"

# --- zh-cn sheet -------------------------------------------------------
# ebb52c47 row (row 3): new Correspond Handoff Datetime / Handback DateTime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-22 20:47:41"
$wsZhCn.Range("K3").Value = "2016-08-22 20:47:58"

# --- de-de sheet ---------------------------------------------------------
# ebb52c47 row (row 3): new Correspond Handoff Datetime / Handback DateTime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-22 20:47:46"
$wsDeDe.Range("K3").Value = "2016-08-22 20:48:15"
